$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.717.17'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.601.95'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.89'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.75'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.827.36'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.590.82'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.523'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.95'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.691.28'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0742'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '210.50'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.16'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0511'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.293.25'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.599'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.15'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.10%  '
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.20'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.04'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.739.35'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.76'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.102'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0515'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.41'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.27%  '
